$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @(2, 6, 1.52),
    @(2, 7, 1.68),
    @(2, 8, 6.8),
    @(2, 9, 9.800000000000001),
    @(2, 11, 4.6),
    @(2, 12, 1.41),
    @(2, 13, 1.07),
    @(2, 14, 3.35),
    @(2, 15, 1.34),
    @(2, 16, 1.8),
    @(2, 17, 2.02),
    @(2, 18, 1.3),
    @(2, 19, 3.25),
    @(2, 20, 2.04),
    @(2, 21, 1.75),
    @(2, 22, 1.13),
    @(2, 24, 16.5),
    @(2, 26, 80),
    @(2, 28, 8.6),
    @(2, 29, 11.5),
    @(2, 31, 160),
    @(2, 32, 10.5),
    @(2, 33, 12.5),
    @(2, 35, 150),
    @(2, 36, 18),
    @(2, 37, 23),
    @(2, 38, 55),
    @(2, 39, 210),
    @(2, 40, 12.5),
    @(2, 41, 250),
    @(3, 6, 2.88),
    @(3, 15, 1.31),
    @(3, 22, 1.59),
    @(3, 23, 1.44),
    @(3, 25, 11.5),
    @(3, 26, 20),
    @(3, 29, 9.4),
    @(3, 33, 16),
    @(3, 34, 20),
    @(4, 6, 1.84),
    @(4, 7, 1.85),
    @(4, 12, 1.3),
    @(4, 16, 2.54),
    @(4, 17, 1.61),
    @(4, 18, 1.62),
    @(4, 21, 2.54),
    @(4, 23, 2.16),
    @(4, 25, 23),
    @(4, 38, 25),
    @(5, 8, 1.29),
    @(5, 26, 9.800000000000001),
    @(5, 27, 12.5),
    @(5, 28, 44),
    @(5, 31, 15),
    @(6, 6, 1.81),
    @(6, 7, 1.91),
    @(6, 8, 4.1),
    @(6, 23, 2.08),
    @(7, 8, 9.199999999999999),
    @(7, 20, 1.78),
    @(7, 21, 2.22),
    @(7, 23, 3.55),
    @(8, 12, 1.33),
    @(8, 16, 2.34),
    @(8, 19, 2.74),
    @(8, 31, 19.5),
    @(8, 35, 29),
    @(8, 39, 70),
    @(8, 40, 29),
    @(8, 41, 12),
    @(9, 10, 4),
    @(9, 11, 4.1),
    @(9, 12, 1.4),
    @(9, 14, 3.95),
    @(9, 15, 1.31),
    @(9, 16, 2.02),
    @(9, 17, 1.93),
    @(9, 18, 1.4),
    @(9, 19, 3.4),
    @(9, 21, 2),
    @(9, 24, 15),
    @(9, 31, 85),
    @(9, 34, 21),
    @(9, 35, 85),
    @(9, 39, 120),
    @(9, 40, 9.800000000000001),
    @(9, 41, 100),
    @(10, 6, 3.2),
    @(10, 7, 4.2),
    @(10, 8, 2.08),
    @(10, 9, 2.32),
    @(10, 10, 3.6),
    @(10, 11, 4.2),
    @(10, 13, 1.05),
    @(10, 14, 4.3),
    @(10, 15, 1.24),
    @(10, 16, 2.14),
    @(10, 17, 1.72),
    @(10, 18, 1.45),
    @(10, 19, 2.78),
    @(10, 20, 1.66),
    @(10, 21, 2.3),
    @(10, 22, 1.75),
    @(10, 23, 1.31),
    @(10, 24, 21),
    @(10, 25, 14.5),
    @(10, 26, 19),
    @(10, 27, 34),
    @(10, 28, 18.5),
    @(10, 29, 11),
    @(10, 30, 13.5),
    @(10, 31, 27),
    @(10, 32, 29),
    @(10, 33, 16.5),
    @(10, 34, 19),
    @(10, 35, 38),
    @(10, 36, 60),
    @(10, 37, 40),
    @(10, 38, 46),
    @(10, 40, 30),
    @(10, 41, 16.5),
    @(11, 6, 1.26),
    @(11, 7, 1.37),
    @(11, 8, 9.4),
    @(11, 9, 14.5),
    @(11, 10, 6),
    @(11, 11, 7.4),
    @(11, 12, 1.19),
    @(11, 14, 6.4),
    @(11, 15, 1.12),
    @(11, 16, 2.88),
    @(11, 17, 1.45),
    @(11, 18, 1.76),
    @(11, 19, 2.02),
    @(11, 20, 1.85),
    @(11, 21, 1.98),
    @(11, 22, 1.07),
    @(11, 23, 3.7),
    @(11, 24, 44),
    @(11, 25, 60),
    @(11, 28, 15),
    @(11, 29, 19),
    @(11, 30, 50),
    @(11, 32, 12),
    @(11, 33, 14),
    @(11, 34, 34),
    @(11, 36, 13),
    @(11, 37, 16.5),
    @(11, 38, 34),
    @(11, 40, 3.9),
    @(12, 6, 1.74),
    @(12, 7, 1.95),
    @(12, 8, 4.3),
    @(12, 9, 6.2),
    @(12, 10, 3.5),
    @(12, 11, 4),
    @(12, 12, 1.33),
    @(12, 13, 1.08),
    @(12, 14, 3.3),
    @(12, 15, 1.34),
    @(12, 16, 1.79),
    @(12, 17, 2),
    @(12, 18, 1.3),
    @(12, 19, 3.65),
    @(12, 20, 1.89),
    @(12, 21, 1.9),
    @(12, 22, 1.2),
    @(12, 23, 2.06),
    @(12, 28, 9.4),
    @(12, 29, 9.800000000000001),
    @(13, 7, 2.12),
    @(13, 18, 1.35),
    @(13, 23, 1.89),
    @(13, 24, 15),
    @(13, 31, 60),
    @(13, 35, 65),
    @(14, 6, 2.64),
    @(14, 7, 2.76),
    @(14, 8, 2.62),
    @(14, 9, 2.74),
    @(14, 11, 4),
    @(14, 15, 1.22),
    @(14, 17, 1.69),
    @(14, 18, 1.52),
    @(14, 19, 2.66),
    @(14, 20, 1.59),
    @(14, 21, 2.48),
    @(14, 22, 1.57),
    @(14, 24, 24),
    @(14, 25, 15.5),
    @(14, 26, 21),
    @(14, 27, 40),
    @(14, 28, 15.5),
    @(14, 30, 13),
    @(14, 31, 27),
    @(14, 32, 22),
    @(14, 33, 13.5),
    @(14, 34, 16),
    @(14, 35, 34),
    @(14, 36, 42),
    @(14, 37, 28),
    @(14, 38, 36),
    @(14, 39, 65),
    @(14, 40, 18.5),
    @(14, 41, 18),
    @(15, 6, 1.65),
    @(15, 7, 1.7),
    @(15, 8, 5.4),
    @(15, 9, 5.5),
    @(15, 11, 4.7),
    @(15, 12, 1.27),
    @(15, 16, 2.34),
    @(15, 17, 1.66),
    @(15, 19, 2.62),
    @(15, 21, 2.2),
    @(15, 23, 2.42),
    @(15, 25, 29),
    @(15, 27, 160),
    @(15, 30, 23),
    @(15, 31, 70),
    @(15, 35, 65),
    @(15, 38, 30),
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}